# New Modified Code for arduino web server
# Appends the latest Arduino weather-station readings to the log sheet:
#   - fills in the missing "Predicted Value" (column P) for rows 25, 34 and 35
#   - completes row 35 (which had been cut short) with its forecast columns
#   - appends three brand-new rows (36, 37, 38) of sensor + forecast data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: fill in the missing Predicted Value ---
$ws.Cells.Item(25, 16).Value = 0

# --- Row 34: fill in the missing Predicted Value ---
$ws.Cells.Item(34, 16).Value = 624.2505119047621

# --- Row 35: finish the row (forecast columns M/N/O were placeholders, plus Predicted Value) ---
$ws.Cells.Item(35, 2).Value = 20.64999961853027
$ws.Cells.Item(35, 3).Value = 82
$ws.Cells.Item(35, 4).Value = 1012
$ws.Cells.Item(35, 5).Value = 4.599999904632568
$ws.Cells.Item(35, 8).Value = 17.29999923706055
$ws.Cells.Item(35, 9).Value = 993.760009765625
$ws.Cells.Item(35, 10).Value = 76
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 1.279999971389771
$ws.Cells.Item(35, 13).Value = 656.39
$ws.Cells.Item(35, 14).Value = 2.41
$ws.Cells.Item(35, 15).Value = 46.6
$ws.Cells.Item(35, 16).Value = 574.6695773809524

# --- Row 36: new reading ---
$ws.Cells.Item(36, 1).Value = "2017.06.03 20.14.05"
$ws.Cells.Item(36, 2).Value = 20.3700008392334
$ws.Cells.Item(36, 3).Value = 100
$ws.Cells.Item(36, 4).Value = 1013
$ws.Cells.Item(36, 5).Value = 1.5
$ws.Cells.Item(36, 6).Value = "2017-06-04T03:00:00"
$ws.Cells.Item(36, 7).Value = "2017-06-04T06:00:00"
$ws.Cells.Item(36, 8).Value = 12.68000030517578
$ws.Cells.Item(36, 9).Value = 994.5700073242188
$ws.Cells.Item(36, 10).Value = 92
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 1.210000038146973
$ws.Cells.Item(36, 13).Value = 166.8
$ws.Cells.Item(36, 14).Value = 24.1
$ws.Cells.Item(36, 15).Value = 459
$ws.Cells.Item(36, 16).Value = 695.5120679320678

# --- Row 37: new reading ---
$ws.Cells.Item(37, 1).Value = "2017.06.03 20.15.06"
$ws.Cells.Item(37, 2).Value = 20.3700008392334
$ws.Cells.Item(37, 3).Value = 100
$ws.Cells.Item(37, 4).Value = 1013
$ws.Cells.Item(37, 5).Value = 1.5
$ws.Cells.Item(37, 6).Value = "2017-06-04T03:00:00"
$ws.Cells.Item(37, 7).Value = "2017-06-04T06:00:00"
$ws.Cells.Item(37, 8).Value = 12.68000030517578
$ws.Cells.Item(37, 9).Value = 994.5700073242188
$ws.Cells.Item(37, 10).Value = 92
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 1.210000038146973
$ws.Cells.Item(37, 13).Value = 166.8
$ws.Cells.Item(37, 14).Value = 24
$ws.Cells.Item(37, 15).Value = 459
$ws.Cells.Item(37, 16).Value = 587.8878928571428

# --- Row 38: new reading (forecast columns not yet populated by the device) ---
$ws.Cells.Item(38, 1).Value = "2017.06.03 20.16.06"
$ws.Cells.Item(38, 2).Value = 20.3700008392334
$ws.Cells.Item(38, 3).Value = 100
$ws.Cells.Item(38, 4).Value = 1013
$ws.Cells.Item(38, 5).Value = 1.5
$ws.Cells.Item(38, 6).Value = "2017-06-04T03:00:00"
$ws.Cells.Item(38, 7).Value = "2017-06-04T06:00:00"
$ws.Cells.Item(38, 8).Value = 12.680000305175781
$ws.Cells.Item(38, 9).Value = 994.5700073242188
$ws.Cells.Item(38, 10).Value = 92
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 1.2100000381469727

# Row 38 has no forecast data yet from the device for M/N/O - leave blank

# --- View state: mirror the selection left behind after the last write ---
$ws.Range("P35").Select() | Out-Null
